$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AllData.trialinfo")

# Insert a new row at row 4 (shifts existing rows 4-17 down to rows 5-18).
$ws.Rows.Item(4).Insert() | Out-Null

# Populate the newly inserted row with the new variable name/explanation
# pair ("ind_trialno" describes the choice-trial count since the last
# emotion induction).
$ws.Cells.Item(4, 1).Value = "ind_trialno"
$ws.Cells.Item(4, 2).Value = "Choice trial number following the most recent emotion induction"

# Leave the workbook with this sheet active and the freshly-typed cell
# selected, matching the author's final view state after the edit.
$ws.Activate() | Out-Null
$ws.Range("B4").Select() | Out-Null
